$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the Talca blueberry sheet.
# It belongs right after the existing 2022-12-12 record (row 114), so
# insert a fresh row at 115 -- this pushes the former rows 115-127 down
# to 116-128 (their data is untouched) and grows the used range to T128.
$ws.Rows.Item(115).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A115").Value2 = 5
$ws.Range("B115").Value2 = "Macroferia Regional de Talca"
$ws.Range("C115").Value2 = "Maule"
$ws.Range("D115").Value2 = 44984
$ws.Range("E115").Value2 = 7
$ws.Range("F115").Value2 = "Fruta"
$ws.Range("G115").Value2 = 100101
$ws.Range("H115").Value2 = "Berries"
$ws.Range("I115").Value2 = 100101001
$ws.Range("J115").Value2 = "Arándano (blue)"
$ws.Range("K115").Value2 = "Sin especificar"
$ws.Range("L115").Value2 = "Primera"
$ws.Range("M115").Value2 = 80
$ws.Range("N115").Value2 = 3000
$ws.Range("O115").Value2 = 3000
$ws.Range("P115").Value2 = 3000
$ws.Range("Q115").Value2 = "$/bandeja 2 kilos"
$ws.Range("R115").Value2 = "Provincia de Curicó"
$ws.Range("S115").Value2 = 1500
$ws.Range("T115").Value2 = 2
